$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update column F (dSF) values for specific rows per repull of data / mean calculation
$ws.Range("F9").Value = 2
$ws.Range("F11").Value = 5
$ws.Range("F18").Value = -3
$ws.Range("F28").Value = 2
$ws.Range("F31").Value = -4
$ws.Range("F36").Value = -10
$ws.Range("F39").Value = -1
